$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("boss_parameter")

$ws.Range("B5").Value = 0.3
$ws.Range("B16").Value = 0.3
$ws.Range("B27").Value = 0.3
$ws.Range("B38").Value = 0.3
$ws.Range("B49").Value = 0.3
$ws.Range("B59").Value = 100
$ws.Range("B60").Value = 0.1
$ws.Range("B63").Value = 200
$ws.Range("B65").Value = 100
$ws.Range("B71").Value = 0.5
$ws.Range("B82").Value = 0.4
$ws.Range("B93").Value = 0.4
$ws.Range("B104").Value = 0.4
$ws.Range("B118").Value = 300
$ws.Range("B119").Value = 3
$ws.Range("B126").Value = 0.5
$ws.Range("B127").Value = 2
$ws.Range("B129").Value = 150
$ws.Range("B137").Value = 0.5
$ws.Range("B140").Value = 150
$ws.Range("B141").Value = 2
$ws.Range("B147").Value = 300
$ws.Range("B148").Value = 0.3
$ws.Range("B149").Value = 4
$ws.Range("B150").Value = 180
$ws.Range("B151").Value = 500
$ws.Range("B159").Value = 0.5
$ws.Range("B162").Value = 300
$ws.Range("B163").Value = 4

# Update selection/view state to match the target workbook
$ws.Range("C120").Select()
